$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 123, pushing the existing 123:223 block down to 129:229
$ws.Rows("123:128").Insert()

# New week of data (Fecha = 2022-12-28, serial 44923) for Sandia / Macroferia Regional de Talca
# Two identical Extra/Primera/Segunda triplets, as per source update.
$newRows = @(
    @{ Row=123; Calidad="Extra";   Vol=3000; Precio=2500 },
    @{ Row=124; Calidad="Primera"; Vol=3000; Precio=2000 },
    @{ Row=125; Calidad="Segunda"; Vol=2000; Precio=1500 },
    @{ Row=126; Calidad="Extra";   Vol=3000; Precio=2500 },
    @{ Row=127; Calidad="Primera"; Vol=3000; Precio=2000 },
    @{ Row=128; Calidad="Segunda"; Vol=2000; Precio=1500 }
)

$fecha = Get-Date -Year 2022 -Month 12 -Day 28 -Hour 0 -Minute 0 -Second 0

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = 5
    $ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($r, 3).Value = "Maule"
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = 7
    $ws.Cells.Item($r, 6).Value = 100112028
    $ws.Cells.Item($r, 7).Value = "Sandia"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = $nr.Calidad
    $ws.Cells.Item($r, 10).Value = $nr.Vol
    $ws.Cells.Item($r, 11).Value = $nr.Precio
    $ws.Cells.Item($r, 12).Value = $nr.Precio
    $ws.Cells.Item($r, 13).Value = $nr.Precio
    $ws.Cells.Item($r, 14).Value = "`$/unidad"
    $ws.Cells.Item($r, 15).Value = "Región del Maule"
    $ws.Cells.Item($r, 16).Value = $nr.Precio
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
